$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (step day +1)
$ws.Range("A1").Value = 45309

# Update D-column prices (steps 1 and 2)
$ws.Range("D32").Value = 2033
$ws.Range("D33").Value = 2430
$ws.Range("D34").Value = 2720
$ws.Range("D35").Value = 2930
$ws.Range("D36").Value = 3100
$ws.Range("D37").Value = 3481
$ws.Range("D38").Value = 3830
$ws.Range("D39").Value = 4100
$ws.Range("D40").Value = 4320
$ws.Range("D46").Value = 2480
$ws.Range("D47").Value = 2770
$ws.Range("D48").Value = 3040
$ws.Range("D49").Value = 3540
$ws.Range("D50").Value = 3851
$ws.Range("D51").Value = 4290
$ws.Range("D52").Value = 4816
$ws.Range("D53").Value = 5020
$ws.Range("D54").Value = 5540
$ws.Range("D55").Value = 6190
$ws.Range("D56").Value = 6721
$ws.Range("D57").Value = 7550
$ws.Range("D58").Value = 8260
$ws.Range("D59").Value = 9050
$ws.Range("D60").Value = 10317.112
$ws.Range("D61").Value = 10949.764
$ws.Range("D67").Value = 4590
$ws.Range("D68").Value = 5235
$ws.Range("D69").Value = 5640
$ws.Range("D70").Value = 5730
$ws.Range("D71").Value = 6470
$ws.Range("D72").Value = 6970
$ws.Range("D73").Value = 7762
$ws.Range("D74").Value = 8335
$ws.Range("D75").Value = 9675
$ws.Range("D76").Value = 10520
$ws.Range("D77").Value = 11710
$ws.Range("D78").Value = 12890
$ws.Range("D79").Value = 14300
$ws.Range("D80").Value = 19730
$ws.Range("D81").Value = 21980
$ws.Range("D87").Value = 8140
$ws.Range("D88").Value = 8400
$ws.Range("D89").Value = 9840
$ws.Range("D90").Value = 10780
$ws.Range("D91").Value = 11900
$ws.Range("D92").Value = 14540
$ws.Range("D93").Value = 15130
$ws.Range("D94").Value = 16800
$ws.Range("D95").Value = 18400
$ws.Range("D96").Value = 20350
